$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Marking" row (row 11): points per right answer, penalty per wrong answer
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Update the "Total" row (row 12): recomputed total marks and the max-marks string
$ws.Range("B12").Value = 92
$ws.Range("E12").Value = "92 / 112"
